# previsao_retorno.xlsx refresh ("atualizei dados bibi e add")
#
# The source data was regenerated a few days later: every "INATIVO - X.X
# meses sem comprar" label advances by 0.1 month, and the two rows whose
# purchase-cadence window rolled over (id_cliente 5985 and 28458) get
# refreshed probability/count/date figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

# Column indexes (header row 1): A=1 id_cliente ... J=10 situacao, K=11 nome
$situacaoCol = 10

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $situacaoCol)
    $val = $cell.Value()
    if ($val -match "^INATIVO - (\d+\.\d+) meses sem comprar$") {
        $meses = [double]$matches[1]
        $novoMeses = $meses + 0.1
        $novoTexto = "INATIVO - {0:N1} meses sem comprar" -f $novoMeses
        $cell.Value = $novoTexto
    }
}

# id_cliente 5985 (LUMA GABRIELLE OLIVEIRA CALDAS) — row 56
$row56 = 56
$ws.Cells.Item($row56, 2).Value = 0.67   # prob_media
$ws.Cells.Item($row56, 4).Value = 0.83   # prob_maxima
$ws.Cells.Item($row56, 5).Value = 38     # total_compras_historico
$ws.Cells.Item($row56, 6).Value = 0.83   # regularidade
$ws.Cells.Item($row56, 8).Value = 45863.74923611111   # ultima_compra
$ws.Cells.Item($row56, 9).Value = 45878.74923611111   # proxima_compra

# id_cliente 28458 (BEMOL S/A) — row 116
$row116 = 116
$ws.Cells.Item($row116, 5).Value = 17020   # total_compras_historico
$ws.Cells.Item($row116, 8).Value = 45863.74502314815  # ultima_compra
$ws.Cells.Item($row116, 9).Value = 45864.74502314815  # proxima_compra
